# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E) for each worker's block of rows is
# re-ordered (newest period first), and the "Valor Mora" (F) figures that
# belong to the first/last row of the second worker's block travel with
# their row, which net out to a simple swap once E is reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Worker 1 (ARELIS ESTHER CASTRO ORTIZ) : rows 16-27, periods 2001..2012 ---
# Previously ascending 2001 -> 2012 top-to-bottom; now descending 2012 -> 2001.
$worker1Periods = @(2012,2011,2010,2009,2008,2007,2006,2005,2004,2003,2002,2001)
for ($i = 0; $i -lt $worker1Periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $worker1Periods[$i].ToString()
}

# --- Worker 2 (GINA PAOLA GARCIA MARTINEZ) : rows 28-36, periods 2203..2211 ---
# Previously ascending 2203 -> 2211 top-to-bottom; now descending 2211 -> 2203.
$worker2Periods = @(2211,2210,2209,2208,2207,2206,2205,2204,2203)
for ($i = 0; $i -lt $worker2Periods.Length; $i++) {
    $row = 28 + $i
    $ws.Range("E$row").Value = $worker2Periods[$i].ToString()
}

# "Valor Mora" values stay attached to their row; after the period reorder
# above, rows 28 and 36 end up with their values swapped versus before.
$ws.Range("F28").Value = 28000
$ws.Range("F36").Value = 36000
